$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Escape-Xml($s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# Inserts a new bold / red (EE0000) "List Paragraph"-styled paragraph
# (without any bullet/numbering) right after the paragraph whose text
# equals $afterText.
function Insert-RedBoldParagraphAfter($doc, $afterText, $newText) {
    $idx = Get-ParaIndexByText $doc $afterText
    if ($idx -eq -1) {
        Write-Host "ERROR: could not find paragraph with text: $afterText"
        return
    }
    $p = $doc.Paragraphs.Item($idx)
    $r = $p.Range
    $r.InsertParagraphAfter()
    $newPara = $doc.Paragraphs.Item($idx + 1)
    $nr = $newPara.Range

    $escaped = Escape-Xml $newText

    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
      '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:b/><w:bCs/><w:color w:val="EE0000"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="EE0000"/></w:rPr><w:t>' + $escaped + '</w:t></w:r></w:p>' + `
      '</w:body></w:document>' + `
      '</pkg:xmlData></pkg:part></pkg:package>'
    $nr.InsertXML($frag)
}

# Replaces the plain-text run content of the paragraph whose text equals
# $paraText, re-inserting it either with or without a leading
# <w:lastRenderedPageBreak/> marker -- everything else about the paragraph
# (pPr / numbering / style) is left untouched since only the text-bearing
# sub-range is targeted.
function Set-LastRenderedPageBreak($doc, $paraText, [bool]$addBreak) {
    $idx = Get-ParaIndexByText $doc $paraText
    if ($idx -eq -1) {
        Write-Host "ERROR: could not find paragraph with text: $paraText"
        return
    }
    $p = $doc.Paragraphs.Item($idx)
    $r = $p.Range
    $textRange = $doc.Range($r.Start, $r.Start + $paraText.Length)

    $escaped = Escape-Xml $paraText
    $runInner = $escaped
    if ($addBreak) {
        $runInner = '<w:lastRenderedPageBreak/><w:t>' + $escaped + '</w:t>'
    } else {
        $runInner = '<w:t>' + $escaped + '</w:t>'
    }

    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
      '<w:p><w:r>' + $runInner + '</w:r></w:p>' + `
      '</w:body></w:document>' + `
      '</pkg:xmlData></pkg:part></pkg:package>'
    $textRange.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# 1) Add the four red/bold reviewer-comment paragraphs
# ---------------------------------------------------------------------------

Insert-RedBoldParagraphAfter $d `
    "Does the code completely and correctly implement the design?" `
    "The code in its current form does implement the intended design for the original project, however this design is not up to the standard that I would want it to be."

Insert-RedBoldParagraphAfter $d `
    "Does the code conform to any pertinent coding standards?" `
    "DRY"

Insert-RedBoldParagraphAfter $d `
    "Is the code clearly and adequately documented with an easy-to-maintain commenting style?" `
    "This codebase has virtually no documentation, and this is something I intend to add as part of my first enhancement for the project."

Insert-RedBoldParagraphAfter $d `
    "Are all comments consistent with the code?" `
    "There are no comments in the code, and the negative effects of this are clear as I did have to spend some time re-learning how this program functions."

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker from the "Loops and
#    Branches" heading paragraph to the "Does the code avoid comparing
#    floating-point numbers for equality?" paragraph.
# ---------------------------------------------------------------------------

Set-LastRenderedPageBreak $d "Loops and Branches" $false
Set-LastRenderedPageBreak $d "Does the code avoid comparing floating-point numbers for equality?" $true

Write-Host "Edits complete"
